$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing the original text formatting,
# e.g. trailing zeros like "1.00" or leading zeros in "0.0000101").
$textCells = @("D5","D6","D7","D8","D9","D10","D11","D13","D16","D21","D23","D25","D26","D27","D30","D31","D32","D34","D35","D36","D37","D38","D41","D43","D44","D46","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values from the source diff.
$ws.Range("D2").Value = "71.495.83"
$ws.Range("D3").Value = "3.880.10"
$ws.Range("E3").Value = "  -2.91%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "605.38"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").Value = "174.47"
$ws.Range("E6").Value = "  +6.55%  "
$ws.Range("D7").Value = "0.669"
$ws.Range("E7").Value = "  -2.76%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.750"
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").Value = "0.177"
$ws.Range("D11").Value = "54.24"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "11.44"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("D14").Value = "4.496.04"
$ws.Range("E14").Value = "  -3.04%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.890.68"
$ws.Range("E15").Value = "  -2.79%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "21.01"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("E18").Value = "  -4.01%  "
$ws.Range("D20").Value = "71.369.98"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").Value = "440.02"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("D23").Value = "94.23"
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("E24").Value = "  -4.03%  "
$ws.Range("D25").Value = "13.89"
$ws.Range("E25").Value = "  -4.10%  "
$ws.Range("D26").Value = "11.77"
$ws.Range("E26").Value = "  +3.50%  "
$ws.Range("D27").Value = "4.05"
$ws.Range("E27").Value = "  -5.94%  "
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("D30").Value = "8.77"
$ws.Range("E30").Value = "  +13.29%  "
$ws.Range("D31").Value = "35.21"
$ws.Range("E31").Value = "  -3.64%  "
$ws.Range("D32").Value = "13.60"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("E33").Value = "  -3.75%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "47.92"
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0000101"
$ws.Range("E35").Value = "  +10.88%  "
$ws.Range("D36").Value = "69.80"
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("D37").Value = "631.11"
$ws.Range("E37").Value = "  -2.84%  "
$ws.Range("D38").Value = "0.436"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "3.33"
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "3.20"
$ws.Range("E43").Value = "  +19.98%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.89"
$ws.Range("E44").Value = "  +8.77%  "
$ws.Range("E45").Value = "  -3.90%  "
$ws.Range("D46").Value = "10.24"
$ws.Range("E46").Value = "  -3.62%  "
$ws.Range("E47").Value = "  -3.90%  "
$ws.Range("D48").Value = "2.92"
$ws.Range("E48").Value = "  -12.88%  "
$ws.Range("D49").Value = "2.913.33"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").Value = "0.000280"
$ws.Range("E50").Value = "  +3.33%  "
$ws.Range("D51").Value = "3.23"
$ws.Range("E51").Value = "  -5.57%  "
